$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths (columns B..N / 2..14)
#    (ColumnWidth is expressed in "characters"; values chosen are the closest
#    achievable to the recorded target raw widths.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 11.3333333333333
$ws.Columns.Item(4).ColumnWidth = 12.8333333333333
$ws.Columns.Item(5).ColumnWidth = 9.16666666666667
$ws.Columns.Item(6).ColumnWidth = 8.33333333333333
$ws.Columns.Item(7).ColumnWidth = 8.33333333333333
$ws.Columns.Item(8).ColumnWidth = 12.1666666666667
$ws.Columns.Item(9).ColumnWidth = 13.3333333333333
$ws.Columns.Item(10).ColumnWidth = 66.8333333333333
$ws.Columns.Item(11).ColumnWidth = 37.1666666666667
$ws.Columns.Item(12).ColumnWidth = 8.33333333333333
$ws.Columns.Item(13).ColumnWidth = 12.1666666666667
$ws.Columns.Item(14).ColumnWidth = 8.33333333333333

# ---------------------------------------------------------------------------
# 2. Header row (row 5): bold the existing boxed / filled header band, and
#    give it a trailing (14th) cell to match the new column N.
# ---------------------------------------------------------------------------
$ws.Range("A5:M5").Font.Bold = $true
$ws.Range("N5").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Data rows 6 & 7 - a new "Xe hop dong" marker column (C) is populated
#    with "X" and every following value shifts one column to the right.
# ---------------------------------------------------------------------------
# Row 6
$ws.Range("C6").Value = "X"
$ws.Range("E6").Value = "Toyota         "
$ws.Range("F6").Value = "Xe con              "
$ws.Range("G6").Value = "2018      "
$ws.Range("H6").Value = "Xe t?i B"
$ws.Range("I6").Value = "B1   "
$ws.Range("J6").Value = "Tạm thời không có trong db"
$ws.Range("K6").Value = "Tạm thời không có trong db"
$ws.Range("M6").Value = "Quang Dat                     "

# Row 7
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "SDK002    "
$ws.Range("C7").Value = "X"
$ws.Range("E7").Value = "Toyota         "
$ws.Range("F7").Value = "Xe  Vip             "
$ws.Range("G7").Value = "2023      "
$ws.Range("H7").Value = "Mô tô"
$ws.Range("I7").Value = "A3   "
$ws.Range("J7").Value = "Tạm thời không có trong db"
$ws.Range("K7").Value = "Tạm thời không có trong db"
$ws.Range("M7").Value = "Ho Quang Dat                  "

# Old column D content (now obsolete / blank - "Xe chu so huu" marker is
# empty for both rows) : clear it out, row 6 drops the cell altogether while
# row 7 keeps an (empty) bordered cell.
$ws.Range("D6").ClearContents()
$ws.Range("D7").ClearContents()

# ---------------------------------------------------------------------------
# 4. Borders: the "X" marker column (C6:C7) gets its own boxed-left/right +
#    hairline-bottom border, matching the rest of the row's hairline grid.
# ---------------------------------------------------------------------------
$ws.Range("C6:C7").Borders.Item(7).LineStyle = 1
$ws.Range("C6:C7").Borders.Item(7).Weight = 2
$ws.Range("C6:C7").Borders.Item(10).LineStyle = 1
$ws.Range("C6:C7").Borders.Item(10).Weight = 2

# Row 7 (last data row) closes the grid with a thin bottom border instead of
# the hairline used on interior rows.
$ws.Range("A7:C7").Borders.Item(9).LineStyle = 1
$ws.Range("A7:C7").Borders.Item(9).Weight = 2
$ws.Range("E7:M7").Borders.Item(9).LineStyle = 1
$ws.Range("E7:M7").Borders.Item(9).Weight = 2

# Blank D7 cell: thin bottom border only (no longer part of the boxed "X"
# column) to align with the rest of row 7's bottom border.
$ws.Range("D7").Borders.Item(9).LineStyle = 1
$ws.Range("D7").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 5. Date stamp cell K8 - refresh the recorded timestamp and give it a
#    dedicated "Ngay dd thang MM nam yyyy" number format, centred.
# ---------------------------------------------------------------------------
$ws.Range("K8").NumberFormat = '"Ngày" dd" tháng "MM" năm "yyyy'
$ws.Range("K8").HorizontalAlignment = -4108
$ws.Range("K8").Value = 45167.8431644444
